# Generate Report for Handback
# Updates the handback-status report with refreshed timestamps/priority
# for the 7f208b2e-... and d2216953-... rows (these two source files
# share the same cached values across the Overview/zh-cn/de-de sheets).

$wb = $excel.ActiveWorkbook

# ---- Sheet "Overview": Latest HO Xliff Generate Date (column G) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-29 02:17:59"   # 7f208b2e-...md row
$wsOverview.Range("G5").Value = "2016-08-29 02:17:59"   # d2216953-...md row

# ---- Sheet "zh-cn" ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column (E): ht -> mt
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
# Correspond Handoff Datetime (column H)
$wsZhCn.Range("H3").Value = "2016-08-29 02:17:54"
$wsZhCn.Range("H5").Value = "2016-08-29 02:17:54"
# Correspond Handback DateTime (column K)
$wsZhCn.Range("K3").Value = "2016-08-29 02:18:14"
$wsZhCn.Range("K5").Value = "2016-08-29 02:18:14"

# ---- Sheet "de-de" ----
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority column (E): ht -> mt
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
# Correspond Handoff Datetime (column H) shares the same cached value as
# the Overview sheet's "Latest HO Xliff Generate Date" for this row.
$wsDeDe.Range("H3").Value = "2016-08-29 02:17:59"
$wsDeDe.Range("H5").Value = "2016-08-29 02:17:59"
# Correspond Handback DateTime (column K)
$wsDeDe.Range("K3").Value = "2016-08-29 02:18:21"
$wsDeDe.Range("K5").Value = "2016-08-29 02:18:21"
